$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header info ---
$ws.Range("C2").Value = "Hartmut"

# B3 holds a long digit string that must stay text (inline/shared string),
# not be auto-converted to a number, and keep its original style (s=8).
$ws.Range("B3").NumberFormat = "@"
$ws.Range("B3").Value = "2570314725427075"
$ws.Range("B2").Copy()
$ws.Range("B3").PasteSpecial(-4122)  # xlPasteFormats - restores original style/number format

$ws.Range("C3").Value = "Mohaupt"

# --- Opening balance date ---
$ws.Range("D5").Value = "KONTOSTAND AM 16.07.2025"

# --- Row 6 ---
$ws.Range("B6").Value = "20.07."
$ws.Range("C6").Value = "21.07."
$ws.Range("D6").Value = "PAYPAL PAACSJ"
$ws.Range("E6").Value = "97,07-"

# --- Row 7 ---
$ws.Range("B7").Value = "24.07."
$ws.Range("C7").Value = "25.07."
$ws.Range("D7").Value = "PAYPAL GOGRJD"
$ws.Range("E7").Value = "26,38-"

# --- Row 8 ---
$ws.Range("B8").Value = "27.07."
$ws.Range("C8").Value = "28.07."
$ws.Range("D8").Value = "RECHNUNG VODAFONE GMBH 79717585"
$ws.Range("E8").Value = "39,86-"

# --- Row 9 ---
$ws.Range("B9").Value = "31.07."
$ws.Range("C9").Value = "01.08."
$ws.Range("D9").Value = "ABSCHLAG STROM Stadtwerke Rosenheim 4620731"
$ws.Range("E9").Value = "87,46-"

# --- Row 10: transaction removed, row becomes blank (matches row 11 styling) ---
$ws.Range("B10").ClearContents()
$ws.Range("C10").ClearContents()
$ws.Range("D10").ClearContents()
$ws.Range("E10").ClearContents()

$ws.Range("E11").Copy()
$ws.Range("E10").PasteSpecial(-4122)  # xlPasteFormats - match the blank row's style (s=12)
$ws.Range("E10").Value = ""

# --- Closing balance ---
$ws.Range("D12").Value = "KONTOSTAND AM 03.08.2025"
$ws.Range("E12").Value = "250,77-"

# --- Next statement date ---
$ws.Range("C13").Value = "IHR NAECHSTER ABRECHNUNGSTERMIN 12.08.2025"
